$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = '"Hey guyss!!! Checking the market trends, Bitcoin is holding steady at $56K. The current trading volume suggests a surge in trust despite the fluctuating market. Our data-driven insights recommend considering investment at the current levels with a stop loss at $54K. Strap in and let''s ride this wave together! 🎢 📈 #cryptocurrency #Bitcoin"'
$ws.Range("C3").Value = 'Hey guyss!!! 🚀 Based on real-time market data, $ETH is currently undervalued. With the increased adoption of Layer 2 solutions, it''s a signal that $ETH might surge soon. My advice? Accumulate below $4K. Keep watching the space! #Crypto #Etherium #InvestSmart👀'
$ws.Range("C4").Value = '"Hey guyss!!! With $BTC trading at a robust $60K and showing strong market dominance, it seems the tables are turned in favor of Bitcoin. Not to mention, some exciting news from the crypto market indicates a potential correlation between current events and BTC rates. Get in while the FOMO''s hot! 🚀📈 #Bitcoin #CryptoNews"'
$ws.Range("C5").Value = '"Hey guyss!!! With $BTC dominance at 60% and bullish institutional inflows, the crypto market looks promising! Look out for the $ETH/BTC ratio at 0.055 which indicates $ETH is undervalued. My analysis suggests an entry point below $4K, with a good chance of it reaching $4.5K soon. Let''s keep riding the crypto wave! 🚀💰💪 #Crypto #Bitcoin #Ethereum"'
$ws.Range("C6").Value = '"Hey guyss!!! $BTC showed an exciting dominance increase to 63%, despite the recent market volatility!🚀 With the constant rise in Hodlers, it''s time to put our focus on Bitcoin! Always remember, careful analysis and strategic action are key to success in crypto!💡🔑 #Bitcoin #CryptoAnalysis"'
